# BC02_HocVienGhiDanh.xlsx — move the "Paid" (Đã đóng) column to the end of
# the table and relabel the totals/paid headers.
#
# Original column layout (row 5 headers / row 6 data-binding row):
#   I = Tổng tiền   (TONGTIEN)
#   J = Đã đóng     (DADONG)
#   K = Miễn giảm (%)
#   L = Miễn giảm (tiền)
#   M = Còn nợ
#
# New layout:
#   I = Học phí     (TONGTIEN, same data, renamed header)
#   J = Miễn giảm (%)
#   K = Miễn giảm (tiền)
#   L = Còn nợ
#   M = Thực đóng   (DADONG, same data, renamed header, moved to the end)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move column J ("Đã đóng" / DADONG) to the end of the table (past column M).
# This shifts K, L, M left by one (into J, K, L) and drops the cut column's
# content into the new column M, with formulas auto-adjusting.
$ws.Columns("J").Cut() | Out-Null
$ws.Columns("N").Insert() | Out-Null

# Relabel the two headers that keep their position/data but change wording.
$ws.Range("I5").Value = "Học phí"
$ws.Range("M5").Value = "Thực đóng"

# The moved "paid" column (and the total/fee column) are right-aligned now
# instead of centered.
$ws.Range("I6:M6").HorizontalAlignment = -4152
$ws.Range("I7:M7").HorizontalAlignment = -4152

# Update the on-screen selection to match the edited state.
$ws.Activate()
$ws.Range("M6").Select()
